# Adds three new data rows to the "Export" sheet of Saldo.xlsx:
#   - 004212476 / MARIA   / 349.25   (just above account 004211807)
#   - 005009992 / ALINE   / 330.17   (just above account 005101676)
#   - 004321016 / JOAQUIM / 0.02     (just below account 004589311 / CLARICE)
#
# Column A ("Conta") holds account numbers with significant leading zeros,
# so they must be written as literal text (not auto-converted to a number).
# We build each text value via a formula (="004212476") and then convert it
# in-place to a plain value with Copy + PasteSpecial xlPasteValues, which
# keeps the cell type as text without leaving a lingering "quote prefix"
# style applied (unlike typing a leading apostrophe into .Value directly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

function Add-SaldoRow {
    param([int]$RowIndex, [string]$Conta, [string]$Nome, [double]$Saldo)

    $ws.Rows.Item($RowIndex).Insert()

    $contaCell = $ws.Cells.Item($RowIndex, 1)
    $contaCell.Formula = '="' + $Conta + '"'
    $contaCell.Copy()
    $contaCell.PasteSpecial(-4163)  # xlPasteValues

    $ws.Cells.Item($RowIndex, 2).Value = $Nome
    $ws.Cells.Item($RowIndex, 3).Value = $Saldo
}

# Insert in top-to-bottom order, accounting for the row-index shift caused
# by each preceding insertion.
Add-SaldoRow 56 "004212476" "MARIA" 349.25
Add-SaldoRow 58 "005009992" "ALINE" 330.17
Add-SaldoRow 199 "004321016" "JOAQUIM" 0.02
